# Generate Report for Handoff
# b.md has been handed off again (new handoff file / timestamp), so its
# status flips from "Handed back: in sync with en-US" to "Ready for
# handoff" on every sheet, and the zh-cn / de-de detail sheets pick up
# the new handoff artifact name + timestamp for that row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row 3 is b.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-22 04:32:51"

# ---------------------------------------------------------------------
# zh-cn sheet - row 3 is b.md
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-22 04:32:48"

# Rebuild the hyperlinks on the zh-cn sheet so the D3 hyperlink's
# display text tracks the new handoff file name (the relationship
# targets themselves are unchanged). Deleting an individual Hyperlinks
# item in place is unreliable, so the whole collection is rebuilt in
# its original order/targets.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d8b1c8fbb2524374e5cf3228bc960e3eaf4bb1cb/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4ece13c8445d98486c362942ec4fe3b50eefff75/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c1f89e241bc5ac48a71ae4d2e7c728b1cd51ec66/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6e5771555255d2b2ae81b7ba1ab7b6abd5da1a34/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d8b1c8fbb2524374e5cf3228bc960e3eaf4bb1cb/e2e/b.md", [Type]::Missing, [Type]::Missing, "b.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4ece13c8445d98486c362942ec4fe3b50eefff75/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c1f89e241bc5ac48a71ae4d2e7c728b1cd51ec66/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6e5771555255d2b2ae81b7ba1ab7b6abd5da1a34/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet - row 3 is b.md
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-22 04:32:51"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d8b1c8fbb2524374e5cf3228bc960e3eaf4bb1cb/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f405f3b00ae9b1d0c688dbe8307df093b3deb750/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/43a9d961ed915bb8ba8c96eed6f87c9026880134/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5363249a95f5a22053d49669b6b102acd93b481a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d8b1c8fbb2524374e5cf3228bc960e3eaf4bb1cb/e2e/b.md", [Type]::Missing, [Type]::Missing, "b.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f405f3b00ae9b1d0c688dbe8307df093b3deb750/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/43a9d961ed915bb8ba8c96eed6f87c9026880134/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5363249a95f5a22053d49669b6b102acd93b481a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

"Report generated for handoff"
